# Merge the three template runs that render the "input date" header cell
# into a single run whose text pipes the date fields through a
# date("dd.MM.YYYY") filter, e.g.:
#   {% if inputDateHeader %}{{ inputDateHeader }}
#   {% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}
#   {% else %}-{% endif %}
# becomes:
#   {% if inputDateHeader %}{{ inputDateHeader | date("dd.MM.YYYY") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date("dd.MM.YYYY") }}){% else %}{% endif %}{% else %}-{% endif %}

$d = $word.ActiveDocument

$old = "{% if inputDateHeader %}{{ inputDateHeader }}{% if paperInputDateHeader %}({{ paperInputDateHeader }}){% else %}{% endif %}{% else %}-{% endif %}"
$new = "{% if inputDateHeader %}{{ inputDateHeader | date(`"dd.MM.YYYY`") }}{% if paperInputDateHeader %} ({{ paperInputDateHeader | date(`"dd.MM.YYYY`") }}){% else %}{% endif %}{% else %}-{% endif %}"

$r = $d.Content
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Assign directly to Range.Text (rather than using Find's Replace
    # parameter) so straight double quotes aren't auto-corrected into
    # curly/smart quotes. After Find.Execute succeeds, $r collapses to
    # the matched range, so this only touches the found text.
    $r.Text = $new
} else {
    Write-Host "WARNING: target text not found"
}
